$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Helper: write a text value into a cell while forcing it to remain
# a TEXT cell (the account numbers have significant leading zeros,
# e.g. "004525587", which must not be coerced into a Number).
# ------------------------------------------------------------------
function Set-TextCell($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $text
}

# ------------------------------------------------------------------
# Helper: locate the row whose column-A value equals $account,
# searching top-down starting at row 2 (row 1 is the header).
# ------------------------------------------------------------------
function Find-AccountRow($account) {
    $used = $ws.UsedRange
    $n = $used.Rows.Count
    for ($i = 2; $i -le $n; $i++) {
        if ($ws.Cells.Item($i, 1).Value() -eq $account) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Row 2: account/name change for the top balance (value stays 85000)
#    004503381 / FREDERICO  ->  004525587 / MARIANA
# ------------------------------------------------------------------
Set-TextCell "A2" "004525587"
$ws.Range("B2").Value = "MARIANA"

# ------------------------------------------------------------------
# 2) MARCIO (008197156 / 50000) row is removed outright, and the
#    JOVINO (005079311 / 9303.57) row right after it is replaced by
#    a new 005009922 / ANA / 6538.79 row (at the same position).
# ------------------------------------------------------------------
$rowMarcio = Find-AccountRow "008197156"
$ws.Rows.Item($rowMarcio).Delete()

$rowJovino = Find-AccountRow "005079311"
Set-TextCell ("A" + $rowJovino) "005009922"
$ws.Range("B" + $rowJovino).Value = "ANA"
$ws.Range("C" + $rowJovino).Value = 6538.79

# ------------------------------------------------------------------
# 3) WAGNER (001804114 / 2195.31) and SILVIA (005547703 / 1272.48)
#    rows are both deleted outright.
# ------------------------------------------------------------------
$rowWagner = Find-AccountRow "001804114"
$ws.Rows.Item($rowWagner).Delete()

$rowSilvia = Find-AccountRow "005547703"
$ws.Rows.Item($rowSilvia).Delete()

# ------------------------------------------------------------------
# 4) The near-the-bottom 005009922 / ANA / 13.81 row is deleted
#    outright (distinct occurrence from the new row added in step 2,
#    which now also carries account 005009922 but sits much higher
#    up the sheet and has a different balance).
# ------------------------------------------------------------------
$used = $ws.UsedRange
$n = $used.Rows.Count
for ($i = $n; $i -ge 2; $i--) {
    if ($ws.Cells.Item($i, 1).Value() -eq "005009922" -and $ws.Cells.Item($i, 3).Value() -eq 13.81) {
        $ws.Rows.Item($i).Delete()
        break
    }
}
